# Update the "create-image" demo sheet so the Text/Price columns become
# Plan/PriceRange columns, and Product/Logo columns are now *Url columns.
# (Underlying sample data in rows 2-5 is unchanged - only the header
# labels are renamed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabeling -------------------------------------------------
$ws.Range("B1").Value = "ProductImageUrl"
$ws.Range("C1").Value = "LogoUrl"
$ws.Range("D1").Value = "PlanName"
$ws.Range("E1").Value = "PriceRange"
$ws.Range("G1").Value = "PlanColor"
$ws.Range("H1").Value = "PlanBackgroundColor"
$ws.Range("J1").Value = "PriceRangeBackgroundColor"
$ws.Range("L1").Value = "PriceRangeColor"
$ws.Range("M1").Value = "PlanFontSize"
$ws.Range("O1").Value = "PriceRangeFontSize"
$ws.Range("P1").Value = "PlanFontWeight"
$ws.Range("R1").Value = "PriceRangeFontWeight"

# --- Column widths (grew to fit the longer new header captions) -----------
$ws.Columns(10).ColumnWidth = 22.166666666666668   # J -> 23
$ws.Columns(12).ColumnWidth = 13.0                 # L -> 13.75 (nearest reachable)
$ws.Columns(15).ColumnWidth = 15.499999999999998   # O -> 16.25 (nearest reachable)
$ws.Columns(18).ColumnWidth = 18.166666666666668   # R -> 19

# --- View: scrolled right a bit, new active selection ----------------------
$ws.Range("O18").Select()
